$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 310; all existing rows 310..436 shift down to 311..437.
$ws.Rows.Item(310).Insert()

# Populate the newly inserted row 310 with the new data point.
# (Same record shape as its neighbours; only Fecha (D) and Volumen (J) are new observations,
# the rest of the fields mirror what used to be row 310 before the shift.)
$ws.Range("A310").Value = 6
$ws.Range("B310").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C310").Value = "Metropolitana"
$ws.Range("D310").Value = 44704
$ws.Range("E310").Value = 13
$ws.Range("F310").Value = 100112039
$ws.Range("G310").Value = "Ciboulette"
$ws.Range("H310").Value = "Sin especificar"
$ws.Range("I310").Value = "Primera"
$ws.Range("J310").Value = 610
$ws.Range("K310").Value = 700
$ws.Range("L310").Value = 800
$ws.Range("M310").Value = 746
$ws.Range("N310").Value = "`$/docena de atados"
$ws.Range("O310").Value = "Región Metropolitana"
$ws.Range("P310").Value = 249
$ws.Range("Q310").Value = 3
$ws.Range("R310").Value = "Hortaliza"
